$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $c = $ws.Range($rangeAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '29.195.80'
$ws.Range('E2').Value = '  +2.59%  '
Set-TextValue 'D3' '1.901.13'
$ws.Range('E3').Value = '  +1.34%  '
Set-TextValue 'D4' '1.005'
$ws.Range('E4').Value = '  -1.58%  '
Set-TextValue 'D5' '315.37'
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('E6').Value = '  -1.59%  '
Set-TextValue 'D7' '0.5116'
$ws.Range('E7').Value = '  -0.03%  '
Set-TextValue 'D8' '0.3929'
$ws.Range('E8').Value = '  -0.79%  '
Set-TextValue 'D9' '0.08423'
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D10' '42.56'
$ws.Range('E10').Value = '  +1.28%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D11' '1.120'
$ws.Range('E11').Value = '  +0.82%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D12' '1.899.22'
$ws.Range('E12').Value = '  +1.27%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D13' '6.237'
$ws.Range('E13').Value = '  -0.18%  '
Set-TextValue 'D14' '20.53'
$ws.Range('E14').Value = '  +0.27%  '
Set-TextValue 'D15' '7.343'
$ws.Range('E15').Value = '  +1.41%  '
Set-TextValue 'D16' '1.005'
$ws.Range('E16').Value = '  -1.62%  '
Set-TextValue 'D17' '93.05'
$ws.Range('E17').Value = '  +2.20%  '
Set-TextValue 'D18' '0.00001109'
$ws.Range('E18').Value = '  -0.20%  '
Set-TextValue 'D19' '0.06727'
$ws.Range('E19').Value = '  -0.71%  '
Set-TextValue 'D20' '17.89'
$ws.Range('E20').Value = '  +0.93%  '
$ws.Range('E21').Value = '  -1.64%  '
Set-TextValue 'D22' '6.035'
$ws.Range('E22').Value = '  +1.57%  '
Set-TextValue 'D23' '29.219.68'
$ws.Range('E23').Value = '  +2.50%  '
Set-TextValue 'D24' '11.15'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('E25').Value = '  -3.31%  '
Set-TextValue 'D26' '2.114.92'
$ws.Range('E26').Value = '  +1.24%  '
$ws.Range('E27').Value = '  -1.03%  '
Set-TextValue 'D28' '20.93'
$ws.Range('E28').Value = '  +0.84%  '
Set-TextValue 'D29' '2.438'
$ws.Range('E29').Value = '  +3.76%  '
Set-TextValue 'D30' '126.44'
$ws.Range('E30').Value = '  -0.45%  '
Set-TextValue 'D31' '1.058'
$ws.Range('E31').Value = '  +1.84%  '
$ws.Range('E32').Value = '  -1.07%  '
Set-TextValue 'D33' '5.914'
$ws.Range('E33').Value = '  +2.66%  '
Set-TextValue 'D34' '3.650'
Set-TextValue 'D35' '0.02473'
$ws.Range('E35').Value = '  +1.57%  '
Set-TextValue 'D36' '0.06600'
$ws.Range('E36').Value = '  +2.01%  '
Set-TextValue 'D37' '9.074'
$ws.Range('E37').Value = '  +2.89%  '
Set-TextValue 'D38' '0.2191'
$ws.Range('E38').Value = '  +0.61%  '
Set-TextValue 'D39' '1.232'
$ws.Range('E39').Value = '  +4.05%  '
Set-TextValue 'D40' '5.094'
$ws.Range('E40').Value = '  +2.14%  '
$ws.Range('E41').Value = '  +1.23%  '
Set-TextValue 'D42' '1.233'
$ws.Range('E42').Value = '  -2.42%  '
Set-TextValue 'D43' '11.27'
$ws.Range('E43').Value = '  +0.45%  '
Set-TextValue 'D44' '1.003'
$ws.Range('E44').Value = '  -1.55%  '
Set-TextValue 'D45' '0.6040'
$ws.Range('E45').Value = '  -0.16%  '
Set-TextValue 'D46' '13.16'
$ws.Range('E46').Value = '  +1.75%  '
Set-TextValue 'D47' '3.682'
$ws.Range('E47').Value = '  -0.86%  '
Set-TextValue 'D48' '2.042'
$ws.Range('E48').Value = '  +2.52%  '
Set-TextValue 'D49' '1.229'
$ws.Range('E49').Value = '  +2.02%  '
Set-TextValue 'D50' '123.11'
$ws.Range('E50').Value = '  +0.80%  '
Set-TextValue 'D51' '1.168'
$ws.Range('E51').Value = '  -3.03%  '
